$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (caseNumber) updates
$ws.Range("R2").Value = 166
$ws.Range("T2").Value = "['01K04873', '023A10838', '026C09059', '026K19034', '036G07553', '039H04732', '03L03830', '043L16294', '044J12261', '106A12593', '106D13687', '106D13999', '106D88900', '106F10599', '106F11867', '106F14206', 'ƒ??']"

# Row 4 (registrationDate) updates
$ws.Range("R4").Value = 116
$ws.Range("T4").Value = "['02/18/2011', '03/28/2012', '04/13/2012', '07/18/2012', '07/26/2010', '1/02/2013', '1/11/2013', '10/23/2012', '10/24/2012', '12/27/2011', '8/09/2012']"

# Row 13 (post) updates
$ws.Range("C13").Value = "int"
$ws.Range("E13").Value = 20
$ws.Range("F13").Value = 10565
$ws.Range("G13").Value = 528.25
$ws.Range("H13").Value = 616.5
$ws.Range("I13").Value = 173.4030519369617
$ws.Range("J13").Value = 30068.61842105263
$ws.Range("K13").Value = 211
$ws.Range("L13").Value = 834
